$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New upload-history rows (33-46) appended below the existing data.
# Text-like columns (dates, zero-padded periods, etc.) are protected with
# a temporary "@" (text) number format while the value is written so Excel
# does not auto-convert them to dates/numbers; the format is then reset back
# to the default "Normal" style to match the rest of the sheet.

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$row = 34
$ws.Cells.Item($row, 1).Value = 33
Set-TextCell $ws $row 2 "2022-11-03"
Set-TextCell $ws $row 3 "MER HEALTH SYSTEM"
Set-TextCell $ws $row 4 "EMR SITE"
Set-TextCell $ws $row 5 "202209"
Set-TextCell $ws $row 6 "FTLV9nOnAFC"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 35
$ws.Cells.Item($row, 1).Value = 34
Set-TextCell $ws $row 2 "2022-11-03"
Set-TextCell $ws $row 3 "MER HEALTH SYSTEM"
Set-TextCell $ws $row 4 "EMR SITE"
Set-TextCell $ws $row 5 "202008"
Set-TextCell $ws $row 6 "FTLV9nOnAFC"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 36
$ws.Cells.Item($row, 1).Value = 35
Set-TextCell $ws $row 2 "2022-11-04"
Set-TextCell $ws $row 3 "MER HEALTH SYSTEM"
Set-TextCell $ws $row 4 "EMR SITE"
Set-TextCell $ws $row 5 "202112"
Set-TextCell $ws $row 6 "FTLV9nOnAFC"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 37
$ws.Cells.Item($row, 1).Value = 36
Set-TextCell $ws $row 2 "2022-11-04"
Set-TextCell $ws $row 3 "MER HEALTH SYSTEM"
Set-TextCell $ws $row 4 "EMR SITE"
Set-TextCell $ws $row 5 "202009"
Set-TextCell $ws $row 6 "FTLV9nOnAFC"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 38
$ws.Cells.Item($row, 1).Value = 37
Set-TextCell $ws $row 2 "2022-11-08"
Set-TextCell $ws $row 3 "MER C&T"
Set-TextCell $ws $row 4 "DSD TX NEW"
Set-TextCell $ws $row 5 "2022Q1"
Set-TextCell $ws $row 6 "KxezVOQ2TVR"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 39
$ws.Cells.Item($row, 1).Value = 38
Set-TextCell $ws $row 2 "2022-11-09"
Set-TextCell $ws $row 3 "MER C&T"
Set-TextCell $ws $row 4 "DSD TX NEW"
Set-TextCell $ws $row 5 "2022Q3"
Set-TextCell $ws $row 6 "DoyPc35A7zI"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 40
$ws.Cells.Item($row, 1).Value = 39
Set-TextCell $ws $row 2 "2022-11-09"
Set-TextCell $ws $row 3 "MER C&T"
Set-TextCell $ws $row 4 "DSD TX NEW"
Set-TextCell $ws $row 5 "2022Q2"
Set-TextCell $ws $row 6 "DoyPc35A7zI"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 41
$ws.Cells.Item($row, 1).Value = 40
Set-TextCell $ws $row 2 "2022-11-09"
Set-TextCell $ws $row 3 "MER ATS"
Set-TextCell $ws $row 4 "DSD HTS TST"
Set-TextCell $ws $row 5 "2022Q1"
Set-TextCell $ws $row 6 "DoyPc35A7zI"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 42
$ws.Cells.Item($row, 1).Value = 41
Set-TextCell $ws $row 2 "2022-11-09"
Set-TextCell $ws $row 3 "MER C&T"
Set-TextCell $ws $row 4 "DSD TX NEW"
Set-TextCell $ws $row 5 "2022Q1"
Set-TextCell $ws $row 6 "XNYN71gD1ps"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 43
$ws.Cells.Item($row, 1).Value = 42
Set-TextCell $ws $row 2 "2022-11-09"
Set-TextCell $ws $row 3 "MER C&T"
Set-TextCell $ws $row 4 "DSD TX NEW"
Set-TextCell $ws $row 5 "2022Q4"
Set-TextCell $ws $row 6 "aywqWn0Qkf8"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 44
$ws.Cells.Item($row, 1).Value = 43
Set-TextCell $ws $row 2 "2022-11-09"
Set-TextCell $ws $row 3 "MER C&T"
Set-TextCell $ws $row 4 "DSD TX NEW"
Set-TextCell $ws $row 5 "2022Q3"
Set-TextCell $ws $row 6 "aywqWn0Qkf8"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 45
$ws.Cells.Item($row, 1).Value = 44
Set-TextCell $ws $row 2 "2022-11-11"
Set-TextCell $ws $row 3 "MER SMI"
Set-TextCell $ws $row 4 "DSD PMTCT STAT"
Set-TextCell $ws $row 5 "2022Q3"
Set-TextCell $ws $row 6 "aywqWn0Qkf8"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 46
$ws.Cells.Item($row, 1).Value = 45
Set-TextCell $ws $row 2 "2022-11-11"
Set-TextCell $ws $row 3 "MER ATS"
Set-TextCell $ws $row 4 "DSD HTS TST"
Set-TextCell $ws $row 5 "2022Q3"
Set-TextCell $ws $row 6 "aywqWn0Qkf8"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"

$row = 47
$ws.Cells.Item($row, 1).Value = 46
Set-TextCell $ws $row 2 "2022-11-11"
Set-TextCell $ws $row 3 "MER C&T"
Set-TextCell $ws $row 4 "DSD TX NEW"
Set-TextCell $ws $row 5 "2022Q3"
Set-TextCell $ws $row 6 "kt468XD802Y"
$ws.Cells.Item($row, 7).Value = 200
Set-TextCell $ws $row 8 "Sucess"
Set-TextCell $ws $row 9 "https://mail.ccsaude.org.mz:5455/api/33/dataValueSets"
